$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6059.7295
$ws.Range("I43").Value = 6396.364
$ws.Range("K43").Value = 6396.364
$ws.Range("M43").Value = -6327.364

$ws.Range("H51").Value = 7163.357
$ws.Range("J51").Value = 7236.091
$ws.Range("L51").Value = 7236.091
$ws.Range("N51").Value = -8204.091

$ws.Range("H70").Value = 1214.2333
$ws.Range("J70").Value = 1093.3462
$ws.Range("L70").Value = 3280.0386
$ws.Range("N70").Value = -3820.0386

$ws.Range("H73").Value = 1214.2333
$ws.Range("J73").Value = 1093.3462
$ws.Range("L73").Value = 3280.0386
$ws.Range("N73").Value = -5152.0386

$ws.Range("H74").Value = 6293.647
$ws.Range("I74").Value = 6570.857
$ws.Range("K74").Value = 6570.857
$ws.Range("M74").Value = -5634.857

$ws.Range("H77").Value = 6293.647
$ws.Range("I77").Value = 6570.857
$ws.Range("K77").Value = 32854.285
$ws.Range("M77").Value = -28174.285

$ws.Range("H88").Value = 2618.8
$ws.Range("J88").Value = 2656.7144
$ws.Range("L88").Value = 2656.7144
$ws.Range("N88").Value = -3468.7144

$ws.Range("H91").Value = 2618.8
$ws.Range("J91").Value = 2656.7144
$ws.Range("L91").Value = 2656.7144
$ws.Range("N91").Value = -5464.7144

$ws.Range("H100").Value = 3134
$ws.Range("I100").Value = 3148.5715
$ws.Range("K100").Value = 3148.5715
$ws.Range("M100").Value = -2607.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3677.9546
$ws.Range("I74").Value = 1206.1052
$ws.Range("J74").Value = 19333
$ws.Range("K74").Value = 1206.1052
$ws.Range("L74").Value = 19333
$ws.Range("M74").Value = -332.1052
$ws.Range("N74").Value = -21081

$ws.Range("H77").Value = 3677.9546
$ws.Range("I77").Value = 1206.1052
$ws.Range("J77").Value = 19333
$ws.Range("K77").Value = 6030.526
$ws.Range("L77").Value = 96665
$ws.Range("M77").Value = -1662.526
$ws.Range("N77").Value = -105401

$ws.Range("H102").Value = 1161.6666
$ws.Range("I102").Value = 1213.2609
$ws.Range("J102").Value = 865
$ws.Range("K102").Value = 1213.2609
$ws.Range("L102").Value = 865
$ws.Range("M102").Value = 408.7391
$ws.Range("N102").Value = -4109

$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 30306016
$ws.Range("I132").Value = 1690.875
$ws.Range("K132").Value = 5072.625
$ws.Range("M132").Value = -2542.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 548.129
$ws.Range("I94").Value = 604.9091
$ws.Range("J94").Value = 409.33334
$ws.Range("K94").Value = 604.9091
$ws.Range("L94").Value = 409.33334
$ws.Range("M94").Value = -153.9091
$ws.Range("N94").Value = -1311.33334

$ws.Range("H99").Value = 1582.5
$ws.Range("I99").Value = 1424.125
$ws.Range("K99").Value = 1424.125
$ws.Range("M99").Value = 73.875

$ws.Range("H105").Value = 2864.8
$ws.Range("I105").Value = 3165.3333
$ws.Range("K105").Value = 3165.3333
$ws.Range("M105").Value = -1418.3333

$ws.Range("H134").Value = 3373
$ws.Range("I134").Value = 1601.3846
$ws.Range("J134").Value = 7979.2
$ws.Range("K134").Value = 4804.1538
$ws.Range("L134").Value = 23937.6
$ws.Range("M134").Value = -2269.1538
$ws.Range("N134").Value = -29007.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 669.9286
$ws.Range("I22").Value = 486.66666
$ws.Range("J22").Value = 999.8
$ws.Range("K22").Value = 486.66666
$ws.Range("L22").Value = 999.8
$ws.Range("M22").Value = -136.66666
$ws.Range("N22").Value = -1699.8

$ws.Range("H31").Value = 5876.4863
$ws.Range("I31").Value = 5399.8887
$ws.Range("K31").Value = 5399.8887
$ws.Range("M31").Value = -5104.8887

$ws.Range("H34").Value = 5876.4863
$ws.Range("I34").Value = 5399.8887
$ws.Range("K34").Value = 5399.8887
$ws.Range("M34").Value = -5197.8887

$ws.Range("H58").Value = 4553.387
$ws.Range("I58").Value = 3879.7083
$ws.Range("J58").Value = 6863.143
$ws.Range("K58").Value = 3879.7083
$ws.Range("L58").Value = 6863.143
$ws.Range("M58").Value = -3676.7083
$ws.Range("N58").Value = -7269.143

$ws.Range("H92").Value = 172308.67
$ws.Range("J92").Value = 172308.67
$ws.Range("L92").Value = 172308.67
$ws.Range("N92").Value = -177300.67

$ws.Range("H107").Value = 1171.138
$ws.Range("I107").Value = 993.3043
$ws.Range("K107").Value = 993.3043
$ws.Range("M107").Value = 926.6957

$ws.Range("H110").Value = 199500
$ws.Range("J110").Value = 199500
$ws.Range("L110").Value = 199500
$ws.Range("N110").Value = -207680

$ws.Range("H136").Value = 4553.387
$ws.Range("I136").Value = 3879.7083
$ws.Range("J136").Value = 6863.143
$ws.Range("K136").Value = 11639.1249
$ws.Range("L136").Value = 20589.429
$ws.Range("M136").Value = -9089.124899999999
$ws.Range("N136").Value = -25689.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 166773.17
$ws.Range("I11").Value = 200107.8
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 600323.3999999999
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -600183.3999999999
$ws.Range("N11").Value = -580

$ws.Range("H12").Value = 270.125
$ws.Range("I12").Value = 66
$ws.Range("K12").Value = 198
$ws.Range("M12").Value = -25

$ws.Range("H103").Value = 455.6
$ws.Range("I103").Value = 114
$ws.Range("K103").Value = 342
$ws.Range("M103").Value = 537

$ws.Range("H113").Value = 1250.0588
$ws.Range("J113").Value = 1250.0588
$ws.Range("L113").Value = 3750.1764
$ws.Range("N113").Value = -8090.1764

$ws.Range("H133").Value = 11309.9
$ws.Range("I133").Value = 5639.8
$ws.Range("K133").Value = 16919.4
$ws.Range("M133").Value = -11859.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 12000
$ws.Range("I21").Value = 12000
$ws.Range("K21").Value = 12000
$ws.Range("M21").Value = -11827

$ws.Range("H30").Value = 12000
$ws.Range("I30").Value = 12000
$ws.Range("K30").Value = 12000
$ws.Range("M30").Value = -11895

$ws.Range("H70").Value = 5597.9473
$ws.Range("I70").Value = 5629
$ws.Range("K70").Value = 5629
$ws.Range("M70").Value = -5359

$ws.Range("H73").Value = 5597.9473
$ws.Range("I73").Value = 5629
$ws.Range("K73").Value = 5629
$ws.Range("M73").Value = -4693

$ws.Range("H102").Value = 3128.36
$ws.Range("I102").Value = 2330.1428
$ws.Range("K102").Value = 2330.1428
$ws.Range("M102").Value = -708.1428000000001

$ws.Range("H136").Value = 14216.5
$ws.Range("J136").Value = 14872.546
$ws.Range("L136").Value = 44617.638
$ws.Range("N136").Value = -49717.638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2232.1667
$ws.Range("J22").Value = 2688.2222
$ws.Range("L22").Value = 2688.2222
$ws.Range("N22").Value = -3278.2222

$ws.Range("H27").Value = 2232.1667
$ws.Range("J27").Value = 2688.2222
$ws.Range("L27").Value = 2688.2222
$ws.Range("N27").Value = -2902.2222

$ws.Range("H40").Value = 2583.6667
$ws.Range("I40").Value = 2098.7334
$ws.Range("K40").Value = 2098.7334
$ws.Range("M40").Value = -1962.7334

$ws.Range("H61").Value = 4385.4
$ws.Range("I61").Value = 1877.5
$ws.Range("K61").Value = 1877.5
$ws.Range("M61").Value = -1675.5

$ws.Range("H82").Value = 1460.091
$ws.Range("J82").Value = 1831
$ws.Range("L82").Value = 1831
$ws.Range("N82").Value = -2553

$ws.Range("H85").Value = 1460.091
$ws.Range("J85").Value = 1831
$ws.Range("L85").Value = 1831
$ws.Range("N85").Value = -4327

$ws.Range("H113").Value = 4385.4
$ws.Range("I113").Value = 1877.5
$ws.Range("K113").Value = 1877.5
$ws.Range("M113").Value = 292.5

$ws.Range("H122").Value = 14000
$ws.Range("J122").Value = 14000
$ws.Range("L122").Value = 42000
$ws.Range("N122").Value = -46900

$ws.Range("H136").Value = 35720996
$ws.Range("I136").Value = 6364.6665
$ws.Range("K136").Value = 19093.9995
$ws.Range("M136").Value = -16543.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2321.1924
$ws.Range("I122").Value = 2235.9565
$ws.Range("K122").Value = 6707.869499999999
$ws.Range("M122").Value = -4257.869499999999
